$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 186, shifting existing rows 186:192 down to 187:193.
$ws.Rows.Item(186).Insert()

# Populate the newly inserted row 186 with the new record.
$ws.Cells.Item(186, 1).Value = 9
$ws.Cells.Item(186, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(186, 3).Value = "Metropolitana"
$ws.Cells.Item(186, 4).Value = 45008
$ws.Cells.Item(186, 5).Value = 13
$ws.Cells.Item(186, 6).Value = "Fruta"
$ws.Cells.Item(186, 7).Value = 100103
$ws.Cells.Item(186, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(186, 9).Value = 100103002
$ws.Cells.Item(186, 10).Value = "Ciruela"
$ws.Cells.Item(186, 11).Value = "Angeleno"
$ws.Cells.Item(186, 12).Value = "Primera"
$ws.Cells.Item(186, 13).Value = 500
$ws.Cells.Item(186, 14).Value = 12000
$ws.Cells.Item(186, 15).Value = 12500
$ws.Cells.Item(186, 16).Value = 12220
$ws.Cells.Item(186, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(186, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(186, 19).Value = 679
$ws.Cells.Item(186, 20).Value = 18
